# Changes of 24th May 2022
# Updates the FedEx shipment tracking numbers (column P) for rows 2-26,
# and the ActualRate/Result (columns Q/R) for row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column P: ShipmentTracking numbers -------------------------------
# These look like plain numbers, so Excel would otherwise silently coerce
# them to numeric cells. Pre-format the range as Text so they are stored
# as shared strings (matching the original file's cell typing), then
# strip the temporary formatting back off once the values are in place
# so the cells end up styleless again, same as before the edit.
$ws.Range("P2:P26").NumberFormat = "@"

$ws.Range("P2").Value = "320018624657"
$ws.Range("P3").Value = "320018621073"
$ws.Range("P4").Value = "320018621100"
$ws.Range("P5").Value = "320018621121"
$ws.Range("P6").Value = "320018621165"
$ws.Range("P7").Value = "320018621187"
$ws.Range("P8").Value = "320018621213"
$ws.Range("P9").Value = "320018621235"
$ws.Range("P10").Value = "320018621268"
$ws.Range("P11").Value = "320018621280"
$ws.Range("P12").Value = "320018621327"
$ws.Range("P13").Value = "320018621349"
$ws.Range("P14").Value = "320018621371"
$ws.Range("P15").Value = "320018621393"
$ws.Range("P16").Value = "320018621420"
$ws.Range("P17").Value = "320018621441"
$ws.Range("P18").Value = "320018621485"
$ws.Range("P19").Value = "320018621500"
$ws.Range("P20").Value = "320018621533"
$ws.Range("P21").Value = "320018621555"
$ws.Range("P22").Value = "320018621588"
$ws.Range("P23").Value = "320018621599"
$ws.Range("P24").Value = "320018621603"
$ws.Range("P25").Value = "320018621614"
$ws.Range("P26").Value = "320018621625"

$ws.Range("P2:P26").ClearFormats()

# --- Row 24: ActualRate (Q24) & Result (R24) ---------------------------
# Q24 looks like a currency amount, so Excel would coerce it to a number
# and drag in a new number format style. Same trick: format as text,
# write it, then fix the per-cell style up by pasting the (styleless)
# format from a sibling Q-column cell so Q24 ends up relying on the
# column's default style again, exactly like its neighbours.
$ws.Range("Q24").NumberFormat = "@"
$ws.Range("Q24").Value = "$248.51"
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("Q24").PasteSpecial(-4122) | Out-Null

$ws.Range("R24").Value = "FAIL"

# --- Dimension widens to column T --------------------------------------
# The saved worksheet's <dimension> grows from A1:S26 to A1:T26 even
# though no data lives in column T. Touching a cell's formatting in
# column T (without leaving any value behind) reproduces that expanded
# used-range while leaving the sheet's visible content untouched.
$ws.Range("T2").Font.Bold = $false
